$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 49: column C (IMG) changes from "x" to "v"
$ws.Cells.Item(49, 3).Value = "v"

# New card rows 56-61 (CARD, DESC(+UPG), IMG)
$ws.Cells.Item(56, 1).Value = "ExecutionStrike"
$ws.Cells.Item(56, 2).Value = "v"
$ws.Cells.Item(56, 3).Value = "x"

$ws.Cells.Item(57, 1).Value = "FirstStrike"
$ws.Cells.Item(57, 2).Value = "v"
$ws.Cells.Item(57, 3).Value = "x"

$ws.Cells.Item(58, 1).Value = "MobileFortress"
$ws.Cells.Item(58, 2).Value = "v"
$ws.Cells.Item(58, 3).Value = "x"

$ws.Cells.Item(59, 1).Value = "Catharsis"
$ws.Cells.Item(59, 2).Value = "v"
$ws.Cells.Item(59, 3).Value = "x"

$ws.Cells.Item(60, 1).Value = "Ambush"
$ws.Cells.Item(60, 2).Value = "v"
$ws.Cells.Item(60, 3).Value = "v"

$ws.Cells.Item(61, 1).Value = "ClumsyStrike"
$ws.Cells.Item(61, 2).Value = "v"
$ws.Cells.Item(61, 3).Value = "x"

# Update active selection to match the authored state
$ws.Range("C61").Select()
